$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J (copy formatting from existing header cell H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-26
$data = @(
    @(6, 8),
    @(8, 8),
    @(8, 9),
    @(8, 9),
    @(9, 9),
    @(5, 6),
    @(9, 9),
    @(6, 8),
    @(7, 7),
    @(8, 8),
    @(7, 8),
    @(8, 8),
    @(7, 8),
    @(7, 8),
    @(7, 8),
    @(10, 10),
    @(7, 7),
    @(5, 8),
    @(7, 7),
    @(7, 7),
    @(6, 7),
    @(5, 5),
    @(7, 7),
    @(5, 5),
    @(4, 4)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
